# Applies the "fixed glitch" edits to the Usability Testing slides
# (slide 9 = "Content Placeholder 2" id=4, slide 10 = "Content Placeholder 2" id=3).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 9 : text-only fixes (no paragraphs added/removed)
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(2)
$tr9 = $sh9.TextFrame.TextRange

# Para 2: "<TAB>Easy/Moderate" -> "Neutral/Agree"
$tr9.Paragraphs(2,1).Text = "Neutral/Agree"

# Para 5: "<TAB>Agree" -> "Agree"
$tr9.Paragraphs(5,1).Text = "Agree"

# Para 8: "I prefer using Smart-Waiter over traditional sense"
#   -> split into two runs with the same combined text:
#      "I prefer using Smart-Waiter over traditional " + "sense"
$para8 = $tr9.Paragraphs(8,1)
$prefix = "I prefer using Smart-Waiter over traditional "
$suffixLen = $para8.Length - $prefix.Length
$tail = $para8.Characters($prefix.Length + 1, $suffixLen)
$tail.Text = "XXXXX"
$tail2 = $para8.Characters($prefix.Length + 1, 5)
$tail2.Text = "sense"

# Para 9: "<TAB>No Preference/Agree" -> "Agree"
$tr9.Paragraphs(9,1).Text = "Agree"

# Para 12: "The interface of the system was pleasant"
#   -> split into two runs with the same combined text:
#      "The interface of the system was " + "pleasant"
$para12 = $tr9.Paragraphs(12,1)
$prefix2 = "The interface of the system was "
$suffixLen2 = $para12.Length - $prefix2.Length
$tail3 = $para12.Characters($prefix2.Length + 1, $suffixLen2)
$tail3.Text = "YYYYYYYY"
$tail4 = $para12.Characters($prefix2.Length + 1, 8)
$tail4.Text = "pleasant"

# Para 13: "<TAB>Disagree" -> "Disagree"
$tr9.Paragraphs(13,1).Text = "Disagree"

# ---------------------------------------------------------------------
# Slide 10 : text fixes + paragraph removals
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(2)
$tr10 = $sh10.TextFrame.TextRange

# Para 2: "<TAB>Neutral/Disagree" -> "Disagree"
$tr10.Paragraphs(2,1).Text = "Disagree"

# Para 5: "<TAB>" + "Agree" -> " " + "      Agree" (keep two runs)
$para5 = $tr10.Paragraphs(5,1)
$c1 = $para5.Characters(1,1)
$c1.Text = " "
$c2 = $para5.Characters(2,5)
$c2.Text = "      Agree"

# Remove the now-duplicated "organization of information" / "Neutral/Disagree"
# question block (paragraphs 7,8,9), and the "Agree/Neutral" trailing blank
# paragraph (paragraph 12) -- delete from the highest index down so earlier
# indices stay valid.
$tr10.Paragraphs(12,1).Delete()
$tr10.Paragraphs(9,1).Delete()
$tr10.Paragraphs(8,1).Delete()
$tr10.Paragraphs(7,1).Delete()

# Former paragraph 11 ("Agree/Neutral") is now paragraph 8 -> "Agree"
$tr10.Paragraphs(8,1).Text = "Agree"
